$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top; this shifts all existing rows (1-49) down to (3-51)
# and automatically keeps data intact (including the style on the former row 1).
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# New row 1: numeric column indexes 0..13
for ($col = 1; $col -le 14; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}

# Give new row 1 the bold / bordered / centered-top header look that the
# original header row had (and that the inserted row carried down to row 3).
$newHeader = $ws.Range("A1:N1")
$newHeader.Font.Bold = $true
$newHeader.Borders.LineStyle = 1
$newHeader.HorizontalAlignment = -4108
$newHeader.VerticalAlignment = -4160

# New row 2: blank, except E2 = "Drive"
$ws.Cells.Item(2, 5).Value = "Drive"

# Row 3 (previously row 1, the original text header) should no longer carry
# the bold/border/centered style - revert it back to the default/normal style.
$ws.Range("A3:N3").Style = "Normal"

Write-Host "Applied header restructuring"
